$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the rolled-value counts (column B) that changed since the last stats update.
$ws.Range("B3").Value = 18
$ws.Range("B4").Value = 15
$ws.Range("B5").Value = 15
$ws.Range("B12").Value = 24
$ws.Range("B15").Value = 21
$ws.Range("B18").Value = 16
$ws.Range("B20").Value = 17

# Restore the view to the top of the sheet with B4 selected.
$ws.Activate()
$ws.Range("B4").Select()
